$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offset = 5/6

$ws.Range("A:A").ColumnWidth = 30 - $offset
$ws.Range("B:H").ColumnWidth = 20 - $offset
$ws.Range("I:I").ColumnWidth = 6 - $offset
